# Support I18N for report
# Replace the hard-coded Japanese header labels in row 1 with I18N
# template expressions, matching the data-row template expressions
# already used in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: column headers -> now driven by I18N lookups instead of literal text
$ws.Range("A1").Value = "&=I18N.Code(bean)"
$ws.Range("B1").Value = "&=I18N.Name(bean)"

# Row 2: data template expressions remain unchanged
$ws.Range("A2").Value = "&=item.code"
$ws.Range("B2").Value = "&=item.name"

# Update the active selection to reflect the author's final cursor position
$ws.Range("C2").Select()
